$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A169").Value = "Propiedad destacada x30;  `$587.706;  IVA: `$111.664;  Total: `$699.370;  26-03-21"
$ws.Range("A170").Value = "Plan Escala;  `$146.926;  IVA: `$27.916;  `$174.842;  26-03-21"
$ws.Range("A171").Value = "Plan Escala;  `$146.926;  IVA: `$27.916;  `$174.842;  26-03-21"
